# Apply changes described by the diff:
#  - Update Metadata sheet (URL, Version, Date, Publisher, Description)
#  - Update Elements sheet (Extension row Definition + Constraint(s), Extension.url Fixed Value)

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet ---
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-status"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"
$wsMeta.Range("B11").Value = "Status of the employee based on one or more code systems. Example codes include HIPAA (HipaaEmployeeStatusCodeSystem), Payer (PayerEmployeeStatusCodeSystem) or customer-specific codes."

# --- Elements sheet ---
# Row 2 = "Extension" element
$wsElem.Range("L2").Value = "Status of the employee based on one or more code systems. Example codes include HIPAA (HipaaEmployeeStatusCodeSystem), Payer (PayerEmployeeStatusCodeSystem) or customer-specific codes."
$wsElem.Range("AI2").Value = ""

# Row 5 = "Extension.url" element
$wsElem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-status"
